$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 30.15035247802734
$ws.Range("C5").Value = 61.42125701904297
